# "Generate Report for Handoff"
#
# - Status message changes from "Handoff transform failed" to "Ready for handoff"
#   (this text is shared by the Overview sheet's B2/C2 cells and the B2 cell on
#   each language sheet, since they all hold the same string).
# - Each language sheet (zh-cn, de-de) gets a handoff record filled in for its
#   first data row: the handoff (.xlf) file name + hyperlink, the handoff
#   datetime, and the handoff reason flips from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Handoff transform failed"
$newStatus = "Ready for handoff"

# Overview rollup sheet shows the same status text for each language column.
if ($overview.Range("B2").Value2 -eq $oldStatus) { $overview.Range("B2").Value = $newStatus }
if ($overview.Range("C2").Value2 -eq $oldStatus) { $overview.Range("C2").Value = $newStatus }

# zh-cn handoff details for ddbbd4ca-3d83-49a6-9d26-3bb385f83087.md
$zhcn.Range("B2").Value = $newStatus

$zhName = "ddbbd4ca-3d83-49a6-9d26-3bb385f83087.1daac7d68c43f0cc3440af16bab041e4837f2139.zh-cn.xlf"
$zhUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/c3dd77abf733647f0bf1915376218c26f11ffc6a/e2e/" + $zhName
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), $zhUrl, "", "", $zhName)

$zhcn.Range("D2").Value = "2016-01-28 05:49:45"
$zhcn.Range("H2").Value = "Include"

# de-de handoff details for ddbbd4ca-3d83-49a6-9d26-3bb385f83087.md
$dede.Range("B2").Value = $newStatus

$deName = "ddbbd4ca-3d83-49a6-9d26-3bb385f83087.1daac7d68c43f0cc3440af16bab041e4837f2139.de-de.xlf"
$deUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/c3dd77abf733647f0bf1915376218c26f11ffc6a/e2e/" + $deName
$dede.Hyperlinks.Add($dede.Range("C2"), $deUrl, "", "", $deName)

$dede.Range("D2").Value = "2016-01-28 05:49:56"
$dede.Range("H2").Value = "Include"
